# Generate Report for handoff
#
# A new source file "89d1f5d9-31af-4988-a549-68f6c7bcaf13.md" has become
# ready for handoff. Insert a row for it (immediately above the existing
# ".localization-config" row, which is pushed down by one row) on every
# sheet: the "Overview" summary sheet plus the per-locale "zh-cn" and
# "de-de" detail sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()

$ws.Range("A3").Value = "89d1f5d9-31af-4988-a549-68f6c7bcaf13.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("C4").Value = "Not to be localized"

$ws.Range("A2").Style = "HyperLink"
$ws.Range("A3").Style = "HyperLink"
$ws.Range("A4").Style = "HyperLink"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ff4461a5faa889454f7a8977ee34728bc6d3dd2a/e2e/80044668-9a67-4ea8-bddf-41bd66cd9ed6.md", "", "", "80044668-9a67-4ea8-bddf-41bd66cd9ed6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ff4461a5faa889454f7a8977ee34728bc6d3dd2a/e2e/89d1f5d9-31af-4988-a549-68f6c7bcaf13.md", "", "", "89d1f5d9-31af-4988-a549-68f6c7bcaf13.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ff4461a5faa889454f7a8977ee34728bc6d3dd2a/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()

$ws.Range("A3").Value = "89d1f5d9-31af-4988-a549-68f6c7bcaf13.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "89d1f5d9-31af-4988-a549-68f6c7bcaf13.df8f20b439bed80a06c30fca9716f3feb50ce529.zh-cn.xlf"
$ws.Range("D3").Value = "2016-01-14 02:23:58"
$ws.Range("E3").Value = "80044668-9a67-4ea8-bddf-41bd66cd9ed6.md"
$ws.Range("F3").Value = "80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.zh-cn.xlf"
$ws.Range("G3").Value = "2016-01-14 02:22:51"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Range("A2").Style = "HyperLink"
$ws.Range("C2").Style = "HyperLink"
$ws.Range("E2").Style = "HyperLink"
$ws.Range("F2").Style = "HyperLink"
$ws.Range("A3").Style = "HyperLink"
$ws.Range("C3").Style = "HyperLink"
$ws.Range("E3").Style = "HyperLink"
$ws.Range("F3").Style = "HyperLink"
$ws.Range("A4").Style = "HyperLink"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ff4461a5faa889454f7a8977ee34728bc6d3dd2a/e2e/80044668-9a67-4ea8-bddf-41bd66cd9ed6.md", "", "", "80044668-9a67-4ea8-bddf-41bd66cd9ed6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6c770668245c79c278dbcb9b741046dc5e3e3337/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.zh-cn.xlf", "", "", "80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ba998b7a8f357859094954ca55d15f7e7487c36c/e2e/80044668-9a67-4ea8-bddf-41bd66cd9ed6.md", "", "", "80044668-9a67-4ea8-bddf-41bd66cd9ed6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a50b25a890c03f6967d3dc3c2c040dff6da13974/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.zh-cn.xlf", "", "", "80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ff4461a5faa889454f7a8977ee34728bc6d3dd2a/e2e/89d1f5d9-31af-4988-a549-68f6c7bcaf13.md", "", "", "89d1f5d9-31af-4988-a549-68f6c7bcaf13.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/df8f20b439bed80a06c30fca9716f3feb50ce529/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/89d1f5d9-31af-4988-a549-68f6c7bcaf13.df8f20b439bed80a06c30fca9716f3feb50ce529.zh-cn.xlf", "", "", "89d1f5d9-31af-4988-a549-68f6c7bcaf13.df8f20b439bed80a06c30fca9716f3feb50ce529.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ba998b7a8f357859094954ca55d15f7e7487c36c/e2e/80044668-9a67-4ea8-bddf-41bd66cd9ed6.md", "", "", "80044668-9a67-4ea8-bddf-41bd66cd9ed6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a50b25a890c03f6967d3dc3c2c040dff6da13974/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.zh-cn.xlf", "", "", "80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ff4461a5faa889454f7a8977ee34728bc6d3dd2a/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()

$ws.Range("A3").Value = "89d1f5d9-31af-4988-a549-68f6c7bcaf13.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "89d1f5d9-31af-4988-a549-68f6c7bcaf13.df8f20b439bed80a06c30fca9716f3feb50ce529.de-de.xlf"
$ws.Range("D3").Value = "2016-01-14 02:24:11"
$ws.Range("E3").Value = "80044668-9a67-4ea8-bddf-41bd66cd9ed6.md"
$ws.Range("F3").Value = "80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.de-de.xlf"
$ws.Range("G3").Value = "2016-01-14 02:23:14"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Range("A2").Style = "HyperLink"
$ws.Range("C2").Style = "HyperLink"
$ws.Range("E2").Style = "HyperLink"
$ws.Range("F2").Style = "HyperLink"
$ws.Range("A3").Style = "HyperLink"
$ws.Range("C3").Style = "HyperLink"
$ws.Range("E3").Style = "HyperLink"
$ws.Range("F3").Style = "HyperLink"
$ws.Range("A4").Style = "HyperLink"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ff4461a5faa889454f7a8977ee34728bc6d3dd2a/e2e/80044668-9a67-4ea8-bddf-41bd66cd9ed6.md", "", "", "80044668-9a67-4ea8-bddf-41bd66cd9ed6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7bbb1b353e91750f461af82b44c1d6a6fc92a581/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.de-de.xlf", "", "", "80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c1d61b4d3584144aed466f277dd8e945b13b4e62/e2e/80044668-9a67-4ea8-bddf-41bd66cd9ed6.md", "", "", "80044668-9a67-4ea8-bddf-41bd66cd9ed6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/747e47deec21efc7797cf29382c1a21b449419ca/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.de-de.xlf", "", "", "80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ff4461a5faa889454f7a8977ee34728bc6d3dd2a/e2e/89d1f5d9-31af-4988-a549-68f6c7bcaf13.md", "", "", "89d1f5d9-31af-4988-a549-68f6c7bcaf13.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/df8f20b439bed80a06c30fca9716f3feb50ce529/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/89d1f5d9-31af-4988-a549-68f6c7bcaf13.df8f20b439bed80a06c30fca9716f3feb50ce529.de-de.xlf", "", "", "89d1f5d9-31af-4988-a549-68f6c7bcaf13.df8f20b439bed80a06c30fca9716f3feb50ce529.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c1d61b4d3584144aed466f277dd8e945b13b4e62/e2e/80044668-9a67-4ea8-bddf-41bd66cd9ed6.md", "", "", "80044668-9a67-4ea8-bddf-41bd66cd9ed6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/747e47deec21efc7797cf29382c1a21b449419ca/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.de-de.xlf", "", "", "80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ff4461a5faa889454f7a8977ee34728bc6d3dd2a/.localization-config", "", "", ".localization-config") | Out-Null
